# Weekly update: insert a new daily price record for "Espárragos" at row 4,
# pushing all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 4 (shifts rows 4..34 down to 5..35)
$ws.Rows(4).Insert()

# Populate the new row 4 with the new weekly record
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Macroferia Regional de Talca"
$ws.Range("C4").Value = "Maule"
$ws.Range("D4").Value = 44490
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 300000000
$ws.Range("G4").Value = "Espárragos"
$ws.Range("H4").Value = "Verde"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 850
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = 875
$ws.Range("N4").Value = '$/kilo'
$ws.Range("O4").Value = "Región del Maule"
$ws.Range("P4").Value = 875
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
